$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New transaction rows to append (rows 8-10).
# Columns A (ID) and B (Request ID) are numeric.
# Columns C-G are stored as text in the source data (including the
# numeric-looking Amount/Balance values), so a leading apostrophe is used
# to force Excel to keep them as text rather than auto-converting to numbers.
$newRows = @(
    @{ A = 7; B = 19; C = "100.00"; D = "Credit"; E = "committee@gmail.com"; F = "2023-04-29 01:13:15"; G = "1093.11" },
    @{ A = 8; B = 2;  C = "11.00";  D = "Debit";  E = "faculty1@gmail.com";  F = "2023-04-29 03:23:11"; G = "1082.11" },
    @{ A = 9; B = 22; C = "200.00"; D = "Credit"; E = "committee@gmail.com"; F = "2023-04-29 03:44:20"; G = "1282.11" }
)

$startRow = 8
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $row = $newRows[$i]

    $ws.Cells.Item($r, 1).Value = $row.A
    $ws.Cells.Item($r, 2).Value = $row.B
    $ws.Cells.Item($r, 3).Value = "'" + $row.C
    $ws.Cells.Item($r, 4).Value = $row.D
    $ws.Cells.Item($r, 5).Value = $row.E
    $ws.Cells.Item($r, 6).Value = $row.F
    $ws.Cells.Item($r, 7).Value = "'" + $row.G
}
